$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: replace the volatile CELL()/FIND()/LEFT() array formula with the
# xlSlim WorkbookLocation() UDF (no trailing backslash in its result).
$ws.Range("B1").FormulaArray = "_xll.WorkbookLocation()"

# B2: path separator now supplied explicitly since WorkbookLocation() no
# longer returns a trailing backslash; also no longer volatile (ca removed).
$ws.Range("B2").Formula = "=B1&""\keyword_args.py"""

# B3: selection moves from B17 to B3.
$ws.Range("B3").Select() | Out-Null
